# Update crypto Price (D) and Volume(1h) (E) columns for rows 2-51
# as refreshed by the scheduled GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column values are stored as literal text (e.g. thousands are
# dot-separated, trailing zeros are significant), so numeric-looking
# strings are entered with a leading apostrophe to keep them text
# instead of being auto-converted to numbers by Excel.
$ws.Range("D2").Value = "'69.260.81"
$ws.Range("E2").Value = '  +1.45%  '
$ws.Range("D3").Value = "'3.902.85"
$ws.Range("E3").Value = '  -0.36%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = "'526.11"
$ws.Range("E5").Value = '  +8.19%  '
$ws.Range("D6").Value = "'143.57"
$ws.Range("E6").Value = '  -1.62%  '
$ws.Range("D7").Value = "'0.610"
$ws.Range("E7").Value = '  -2.03%  '
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("D9").Value = "'0.721"
$ws.Range("E9").Value = '  -1.65%  '
$ws.Range("E10").Value = '  +0.99%  '
$ws.Range("D11").Value = "'0.0000331"
$ws.Range("E11").Value = '  -5.23%  '
$ws.Range("D12").Value = "'41.99"
$ws.Range("E12").Value = '  -2.57%  '
$ws.Range("D13").Value = "'4.533.48"
$ws.Range("E13").Value = '  -0.28%  '
$ws.Range("D14").Value = "'10.19"
$ws.Range("E14").Value = '  -4.83%  '
$ws.Range("D15").Value = "'3.915.70"
$ws.Range("E15").Value = '  +0.12%  '
$ws.Range("E16").Value = '  -0.45%  '
$ws.Range("D17").Value = "'1.21"
$ws.Range("E17").Value = '  +6.72%  '
$ws.Range("D18").Value = "'13.73"
$ws.Range("E18").Value = '  -4.16%  '
$ws.Range("D19").Value = "'19.68"
$ws.Range("E19").Value = '  -2.04%  '
$ws.Range("D20").Value = "'69.167.87"
$ws.Range("E20").Value = '  +1.19%  '
$ws.Range("D21").Value = "'427.46"
$ws.Range("E21").Value = '  -1.01%  '
$ws.Range("D22").Value = "'3.32"
$ws.Range("E22").Value = '  -5.71%  '
$ws.Range("D23").Value = "'14.16"
$ws.Range("E23").Value = '  -6.31%  '
$ws.Range("D24").Value = "'87.66"
$ws.Range("E24").Value = '  -0.94%  '
$ws.Range("E25").Value = '  +8.48%  '
$ws.Range("D26").Value = "'11.47"
$ws.Range("E26").Value = '  -2.22%  '
$ws.Range("D27").Value = "'10.59"
$ws.Range("E27").Value = '  -5.75%  '
$ws.Range("D28").Value = "'35.93"
$ws.Range("E28").Value = '  -5.01%  '
$ws.Range("D29").Value = "'693.22"
$ws.Range("E29").Value = '  -3.61%  '
$ws.Range("D30").Value = "'13.07"
$ws.Range("E30").Value = '  -5.11%  '
$ws.Range("E31").Value = '  -4.72%  '
$ws.Range("D32").Value = "'2.80"
$ws.Range("E32").Value = '  -4.28%  '
$ws.Range("D33").Value = "'68.34"
$ws.Range("E33").Value = '  +12.26%  '
$ws.Range("E34").Value = '  +12.64%  '
$ws.Range("D35").Value = "'5.88"
$ws.Range("E35").Value = '  -4.74%  '
$ws.Range("D36").Value = "'39.99"
$ws.Range("E36").Value = '  -4.29%  '
$ws.Range("D37").Value = "'0.0₃0833"
$ws.Range("E37").Value = '  -8.83%  '
$ws.Range("D38").Value = "'0.998"
$ws.Range("E38").Value = '  -0.05%  '
$ws.Range("D39").Value = "'0.146"
$ws.Range("E39").Value = '  +0.96%  '
$ws.Range("E40").Value = '  -0.06%  '
$ws.Range("D41").Value = "'0.0477"
$ws.Range("E41").Value = '  -2.86%  '
$ws.Range("E42").Value = '  +1.74%  '
$ws.Range("D43").Value = "'2.74"
$ws.Range("E43").Value = '  -8.92%  '
$ws.Range("D44").Value = "'2.93"
$ws.Range("E44").Value = '  -6.24%  '
$ws.Range("D45").Value = "'3.35"
$ws.Range("E45").Value = '  +0.10%  '
$ws.Range("E46").Value = '  -1.71%  '
$ws.Range("D47").Value = "'3.03"
$ws.Range("E47").Value = '  +7.73%  '
$ws.Range("D48").Value = "'3.27"
$ws.Range("E48").Value = '  -4.57%  '
$ws.Range("D49").Value = "'142.71"
$ws.Range("E49").Value = '  -1.54%  '
$ws.Range("D50").Value = "'2.04"
$ws.Range("E50").Value = '  -4.72%  '
$ws.Range("D51").Value = "'25.72"
$ws.Range("E51").Value = '  +1.68%  '
